# Applies the "fixed issue of yasumi day" edit:
#  - Slides 10-14 (回路/換える/代替/携わる/連携) are removed entirely.
#  - Slides 1-9 get new kanji/reading/definition text, and the page
#    reference footer changes from "69-70" to "67-68".

$p = $ppt.ActivePresentation

# --- New content for slides 1 through 9 -------------------------------
# Each entry: kanji, reading, definition
$newContent = @(
    @("殺す", "ころす", "to kill, to slay, to murder, to slaughter | to suppress, to block, to hamper, to destroy (e.g. talent), to eliminate (e.g..."),
    @("殺人", "さつじん", "murder, homicide, manslaughter..."),
    @("農薬", "のうやく", "agricultural chemical (i.e. pesticide, herbicide, fungicide, etc.), agrochemical, agrichemical..."),
    @("収入印紙", "しゅうにゅういんし", "revenue stamp..."),
    @("収穫", "しゅうかく", "harvest, crop, ingathering | fruits (of one's labors), gain, result, returns..."),
    @("少量", "しょうりょう", "small quantity, small amount | narrowmindedness..."),
    @("完全", "かんぜん", "perfect, complete..."),
    @("原因", "げんいん", "cause, origin, source..."),
    @("一環", "いっかん", "link (e.g. in a chain of events), part (of a plan, campaign, activities, etc.) | monocyclic...")
)

for ($i = 0; $i -lt $newContent.Count; $i++) {
    $slideIndex = $i + 1
    $s = $p.Slides.Item($slideIndex)

    $word = $newContent[$i][0]
    $reading = $newContent[$i][1]
    $definition = $newContent[$i][2]

    # Shape 1 ("Text 0") holds the headline kanji.
    $titleRange = $s.Shapes.Item(1).TextFrame.TextRange
    $titleRange.Text = ""
    $titleRange.Text = $word

    # Shape 2 ("Text 1") has 2 leading blank paragraphs, reading is
    # in the 3rd paragraph - update only that paragraph in place.
    # (Clearing before assigning avoids the host occasionally
    # fragmenting the new text into multiple runs.)
    $readingPara = $s.Shapes.Item(2).TextFrame.TextRange.Paragraphs(3)
    $readingPara.Text = ""
    $readingPara.Text = $reading

    # Shape 3 ("Text 2") holds the English definition.
    $defRange = $s.Shapes.Item(3).TextFrame.TextRange
    $defRange.Text = ""
    $defRange.Text = $definition

    # Shape 4 ("Text 3") holds the page reference footer.
    $pageRange = $s.Shapes.Item(4).TextFrame.TextRange
    $pageRange.Text = ""
    $pageRange.Text = "67-68"
}

# --- Remove the trailing slides (10-14) --------------------------------
for ($idx = $p.Slides.Count; $idx -ge 10; $idx--) {
    $p.Slides.Item($idx).Delete()
}
